$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 39228.777
$ws.Range("I21").Value = 80019
$ws.Range("J21").Value = 34130
$ws.Range("K21").Value = 80019
$ws.Range("L21").Value = 34130
$ws.Range("M21").Value = -79551
$ws.Range("N21").Value = -35066

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 39228.777
$ws.Range("I23").Value = 80019
$ws.Range("J23").Value = 34130
$ws.Range("K23").Value = 80019
$ws.Range("L23").Value = 34130
$ws.Range("M23").Value = -79785
$ws.Range("N23").Value = -34598

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 4943.95
$ws.Range("I58").Value = 223.4
$ws.Range("J58").Value = 9664.5
$ws.Range("K58").Value = 670.2
$ws.Range("L58").Value = 28993.5
$ws.Range("M58").Value = -520.2
$ws.Range("N58").Value = -29293.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2177370.8
$ws.Range("I74").Value = 3033735.5
$ws.Range("K74").Value = 3033735.5
$ws.Range("M74").Value = -3032799.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 2177370.8
$ws.Range("I77").Value = 3033735.5
$ws.Range("K77").Value = 15168677.5
$ws.Range("M77").Value = -15163997.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 73235.86
$ws.Range("J137").Value = 1946.1538
$ws.Range("L137").Value = 5838.4614
$ws.Range("N137").Value = -10938.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 64223.375
$ws.Range("I74").Value = 81787.67999999999
$ws.Range("K74").Value = 81787.67999999999
$ws.Range("M74").Value = -80913.67999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 64223.375
$ws.Range("I77").Value = 81787.67999999999
$ws.Range("K77").Value = 408938.4
$ws.Range("M77").Value = -404570.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 34751
$ws.Range("I88").Value = 1701.2
$ws.Range("J88").Value = 200000
$ws.Range("K88").Value = 1701.2
$ws.Range("L88").Value = 200000
$ws.Range("M88").Value = -1295.2
$ws.Range("N88").Value = -200812

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 34751
$ws.Range("I91").Value = 1701.2
$ws.Range("J91").Value = 200000
$ws.Range("K91").Value = 1701.2
$ws.Range("L91").Value = 200000
$ws.Range("M91").Value = -297.2
$ws.Range("N91").Value = -202808

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 11389379
$ws.Range("I132").Value = 16127750
$ws.Range("J132").Value = 1548145.9
$ws.Range("K132").Value = 48383250
$ws.Range("L132").Value = 4644437.699999999
$ws.Range("M132").Value = -48380720
$ws.Range("N132").Value = -4649497.699999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52372

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -161856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1789.7142
$ws.Range("I132").Value = 1310.3334
$ws.Range("K132").Value = 3931.0002
$ws.Range("M132").Value = -1401.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 250.73685
$ws.Range("I40").Value = 185.86667
$ws.Range("J40").Value = 494
$ws.Range("K40").Value = 743.46668
$ws.Range("L40").Value = 1976
$ws.Range("M40").Value = -674.46668
$ws.Range("N40").Value = -2114

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1939.2222
$ws.Range("J46").Value = 2667.8333
$ws.Range("L46").Value = 8003.499899999999
$ws.Range("N46").Value = -8185.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 6412213.5
$ws.Range("J58").Value = 7694456.5
$ws.Range("L58").Value = 23083369.5
$ws.Range("N58").Value = -23083625.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2758195
$ws.Range("I64").Value = 1168.6666
$ws.Range("J64").Value = 6066627
$ws.Range("K64").Value = 3505.9998
$ws.Range("L64").Value = 18199881
$ws.Range("M64").Value = -3235.9998
$ws.Range("N64").Value = -18200421

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 2758195
$ws.Range("I67").Value = 1168.6666
$ws.Range("J67").Value = 6066627
$ws.Range("K67").Value = 3505.9998
$ws.Range("L67").Value = 18199881
$ws.Range("M67").Value = -2569.9998
$ws.Range("N67").Value = -18201753

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 166668740
$ws.Range("I82").Value = 796.6667
$ws.Range("J82").Value = 333336670
$ws.Range("K82").Value = 2390.0001
$ws.Range("L82").Value = 1000010010
$ws.Range("M82").Value = -1984.0001
$ws.Range("N82").Value = -1000010822

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 166668740
$ws.Range("I85").Value = 796.6667
$ws.Range("J85").Value = 333336670
$ws.Range("K85").Value = 2390.0001
$ws.Range("L85").Value = 1000010010
$ws.Range("M85").Value = -986.0001000000002
$ws.Range("N85").Value = -1000012818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 6069.846
$ws.Range("J97").Value = 6069.846
$ws.Range("L97").Value = 18209.538
$ws.Range("N97").Value = -19201.538

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 103
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 7630138.5
$ws.Range("I14").Value = 10166851
$ws.Range("J14").Value = 20000
$ws.Range("K14").Value = 10166851
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = -10166683
$ws.Range("N14").Value = -20336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 83339.336
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 83339.336
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 83339.336
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -84397.336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 52337
$ws.Range("J27").Value = 52337
$ws.Range("L27").Value = 52337
$ws.Range("N27").Value = -52669

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1253411.8
$ws.Range("I136").Value = 2503326
$ws.Range("J136").Value = 3497.5
$ws.Range("K136").Value = 7509978
$ws.Range("L136").Value = 10492.5
$ws.Range("M136").Value = -7507428
$ws.Range("N136").Value = -15592.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 15000
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 20000
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 20000
$ws.Range("M76").Value = -4685
$ws.Range("N76").Value = -20630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H79").Value = 15000
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 20000
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 20000
$ws.Range("M79").Value = -3908
$ws.Range("N79").Value = -22184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
